$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 - "Dashboard Overview (1)" title: bump font size to 40pt, extend
# the run text, and re-flow the (shrunk) title placeholder box.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(2)
$shp10.Top = 27.74578857421875
$shp10.Height = 55.73909378051758

$tr10 = $shp10.TextFrame.TextRange
$tr10.Font.Size = 40
$tr10.Runs(3).Text = "iew (1) – Headcount & Distribution"

# ---------------------------------------------------------------------------
# Slide 11 - "Dashboard Overview (2)" title: same treatment.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)
$shp11.Top = 27.745708465576172
$shp11.Height = 55.73909378051758

$tr11 = $shp11.TextFrame.TextRange
$tr11.Font.Size = 40
$tr11.Runs(3).Text = "iew (2) – Salary & Performance"

# ---------------------------------------------------------------------------
# Slide 9 - "Conclusion" body copy: resize/move the content placeholder and
# demote the summary bullet paragraphs to the second outline level.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(1)
$shp9.Top = 153.36004638671875
$shp9.Height = 386.6400451660156

$tr9 = $shp9.TextFrame.TextRange
for ($i = 2; $i -le 8; $i++) {
    $tr9.Paragraphs($i).IndentLevel = 2
}
